# Apply cryptos.xlsx price/volume update
# Commit message: "Updated cryptos list on Sun Jun  2 11:20:31 UTC 2024 with GitHub Actions"
#
# Column D (Price) cells are stored as plain text in the source workbook (inline
# strings, e.g. "67.512.47" uses dots as thousands separators). Excel's COM layer
# auto-converts plain numeric-looking strings into real numbers (dropping trailing
# zeros / switching to scientific notation for tiny values), so those assignments
# are prefixed with a leading apostrophe to force a literal text entry, exactly as
# the original file stores them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''67.545.83'
$ws.Range('D3').Value = '''3.778.96'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('D5').Value = '''597.89'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').Value = '''164.32'
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('E9').Value = '  -1.08%  '
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('D11').Value = '''6.40'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('D12').Value = '''0.0000247'
$ws.Range('E12').Value = '  -1.97%  '
$ws.Range('D13').Value = '''35.43'
$ws.Range('E13').Value = '  -1.72%  '
$ws.Range('D14').Value = '''4.414.21'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').Value = '''3.768.35'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').Value = '''67.574.18'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '''18.30'
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('E18').Value = '  +1.70%  '
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('D20').Value = '''459.20'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('E21').Value = '  -2.94%  '
$ws.Range('D22').Value = '''0.692'
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('D23').Value = '''0.0000145'
$ws.Range('E23').Value = '  -5.10%  '
$ws.Range('D24').Value = '''82.47'
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').Value = '''11.96'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('E26').Value = '  -0.96%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = '''9.91'
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('D29').Value = '''3.927.46'
$ws.Range('D30').Value = '''7.41'
$ws.Range('E30').Value = '  +2.59%  '
$ws.Range('E31').Value = '  -6.07%  '
$ws.Range('E32').Value = '  -3.03%  '
$ws.Range('D33').Value = '''29.02'
$ws.Range('E33').Value = '  -1.91%  '
$ws.Range('D34').Value = '''0.998'
$ws.Range('E34').Value = '  -0.46%  '
$ws.Range('D35').Value = '''8.94'
$ws.Range('E35').Value = '  -1.39%  '
$ws.Range('D36').Value = '''0.0987'
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  -3.35%  '
$ws.Range('D39').Value = '''0.986'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').Value = '''47.38'
$ws.Range('D44').Value = '''43.32'
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('D45').Value = '''0.295'
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('D46').Value = '''151.93'
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('D47').Value = '''8.31'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''27.01'
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '''1.84'
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '''1.34'
$ws.Range('E50').Value = '  +6.33%  '
$ws.Range('D51').Value = '''389.24'
$ws.Range('E51').Value = '  +0.14%  '
